# Insert two new reference/lookup sheets -- "BenefitTypes" and
# "InsuranceTypes" -- right after "ResidenceUses" and before "IncomeTypes".
# This is part of an effort to align all row/column names with the data
# dictionary specs: new BenefitGroup/OfficialBenefitName and
# InsuranceGroup/OfficialInsuranceName lookup tables are added alongside
# the existing IncomeGroup/OfficialIncomeName table.

$wb = $excel.ActiveWorkbook

# Unicode characters used in a couple of the "official" display names.
$rsquo = [char]0x2019   # RIGHT SINGLE QUOTATION MARK (')
$ndash = [char]0x2013   # EN DASH (-)

# --- BenefitTypes sheet ---------------------------------------------------
# Inserted immediately before the current 2nd sheet (IncomeTypes), so it
# lands right after ResidenceUses.
$benefitSheet = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$benefitSheet.Name = "BenefitTypes"

$benefitSheet.Range("A1").Value = "BenefitGroup"
$benefitSheet.Range("B1").Value = "OfficialBenefitName"

$benefitGroups = @("SNAP", "WIC", "TANFChildCare", "TANFTransportation", "OtherTANF", "OtherBenefitsSource")
for ($i = 0; $i -lt $benefitGroups.Count; $i++) {
    $benefitSheet.Cells.Item($i + 2, 1).Value = $benefitGroups[$i]
}

$benefitNames = @(
    "Supplemental Nutrition Assistance Program (SNAP) (Previously known as Food Stamps)",
    "Special Supplemental Nutrition Program for Women, Infants, and Children (WIC)",
    "TANF Child Care Services",
    "TANF Transportation Services",
    "Other TANF-Funded Services",
    "Other Source"
)
for ($i = 0; $i -lt $benefitNames.Count; $i++) {
    $benefitSheet.Cells.Item($i + 2, 2).Value = $benefitNames[$i]
}

# --- InsuranceTypes sheet --------------------------------------------------
# Inserted immediately before the current 3rd sheet (IncomeTypes, now that
# BenefitTypes occupies slot 2).
$insuranceSheet = $wb.Worksheets.Add($wb.Worksheets.Item(3))
$insuranceSheet.Name = "InsuranceTypes"

$insuranceSheet.Range("A1").Value = "InsuranceGroup"
$insuranceSheet.Range("B1").Value = "OfficialInsuranceName"

$insuranceGroups = @("Medicaid", "Medicare", "SCHIP", "VAMedicalServices", "EmployerProvided", "COBRA", "PrivatePay", "StateHealthIns", "IndianHealthServices", "OtherInsurance")
for ($i = 0; $i -lt $insuranceGroups.Count; $i++) {
    $insuranceSheet.Cells.Item($i + 2, 1).Value = $insuranceGroups[$i]
}

$insuranceNames = @(
    "MEDICAID",
    "MEDICARE",
    "State Children" + $rsquo + "s Health Insurance Program",
    "Veteran" + $rsquo + "s Administration (VA) Medical Services",
    "Employer " + $ndash + " Provided Health Insurance",
    "Health Insurance obtained through COBRA",
    "Private Pay Health Insurance",
    "State Health Insurance for Adults",
    "Indian Health Services Program",
    "Other"
)
for ($i = 0; $i -lt $insuranceNames.Count; $i++) {
    $insuranceSheet.Cells.Item($i + 2, 2).Value = $insuranceNames[$i]
}

# InsuranceTypes ends up the active/selected tab, matching the saved
# workbook view (activeTab points at the 3rd sheet, 0-indexed).
$insuranceSheet.Activate()
